# Update gh-pages to output generated at 456a3b4
$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 128
$ws1.Range("F7").Value = 299
$ws1.Range("F9").Value = 2058
$ws1.Range("F10").Value = 359
$ws1.Range("F11").Value = 4928

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 128
$ws4.Range("F9").Value = 299
$ws4.Range("F13").Value = 2058
$ws4.Range("F14").Value = 359
$ws4.Range("F15").Value = 4928
